$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'27.226.50"
$ws.Range("E2").Value = "  +0.15%  "
$ws.Range("D3").Value = "'1.906.03"
$ws.Range("E3").Value = "  +0.12%  "
$ws.Range("E4").Value = "  +0.30%  "
$ws.Range("D5").Value = "'307.70"
$ws.Range("D6").Value = "'1.002"
$ws.Range("E6").Value = "  +0.22%  "
$ws.Range("D7").Value = "'0.5266"
$ws.Range("E7").Value = "  +0.41%  "
$ws.Range("D8").Value = "'0.3821"
$ws.Range("E8").Value = "  +1.45%  "
$ws.Range("D9").Value = "'0.07304"
$ws.Range("E9").Value = "  +0.67%  "
$ws.Range("E10").Value = "  +2.05%  "
$ws.Range("D11").Value = "'0.9053"
$ws.Range("E11").Value = "  +0.46%  "
$ws.Range("D12").Value = "'0.08074"
$ws.Range("E12").Value = "  -4.23%  "
$ws.Range("D13").Value = "'96.02"
$ws.Range("E13").Value = "  +1.24%  "
$ws.Range("D14").Value = "'5.368"
$ws.Range("E14").Value = "  +1.46%  "
$ws.Range("D15").Value = "'1.794.65"
$ws.Range("E15").Value = "  -4.47%  "
$ws.Range("D16").Value = "'1.002"
$ws.Range("E16").Value = "  +0.20%  "
$ws.Range("D17").Value = "'0.000008685"
$ws.Range("E17").Value = "  +0.66%  "
$ws.Range("D18").Value = "'14.74"
$ws.Range("E18").Value = "  +1.03%  "
$ws.Range("E19").Value = "  +0.26%  "
$ws.Range("D20").Value = "'27.264.95"
$ws.Range("E20").Value = "  +0.14%  "
$ws.Range("E21").Value = "  +1.16%  "
$ws.Range("E22").Value = "  +1.94%  "
$ws.Range("D23").Value = "'6.493"
$ws.Range("E23").Value = "  +0.83%  "
$ws.Range("D24").Value = "'2.348"
$ws.Range("E24").Value = "  +3.08%  "
$ws.Range("D25").Value = "'149.84"
$ws.Range("E25").Value = "  +1.73%  "
$ws.Range("D26").Value = "'18.26"
$ws.Range("E26").Value = "  +0.48%  "
$ws.Range("E27").Value = "  -0.56%  "
$ws.Range("D28").Value = "'116.97"
$ws.Range("E28").Value = "  +1.82%  "
$ws.Range("D29").Value = "'4.845"
$ws.Range("E29").Value = "  +0.62%  "
$ws.Range("D30").Value = "'4.882"
$ws.Range("E30").Value = "  -0.67%  "
$ws.Range("D31").Value = "'0.09237"
$ws.Range("E31").Value = "  -0.33%  "
$ws.Range("D32").Value = "'0.8068"
$ws.Range("E32").Value = "  -0.45%  "
$ws.Range("D33").Value = "'0.05066"
$ws.Range("E33").Value = "  -0.02%  "
$ws.Range("D34").Value = "'1.230"
$ws.Range("E34").Value = "  -0.95%  "
$ws.Range("D35").Value = "'2.983"
$ws.Range("E35").Value = "  +1.06%  "
$ws.Range("D36").Value = "'3.380"
$ws.Range("E36").Value = "  -0.39%  "
$ws.Range("D37").Value = "'2.695"
$ws.Range("E37").Value = "  +2.79%  "
$ws.Range("D38").Value = "'0.5738"
$ws.Range("E38").Value = "  -0.23%  "
$ws.Range("D39").Value = "'0.01996"
$ws.Range("E39").Value = "  +0.29%  "
$ws.Range("D40").Value = "'1.086"
$ws.Range("E40").Value = "  +0.94%  "
$ws.Range("D41").Value = "'8.997"
$ws.Range("E41").Value = "  +0.18%  "
$ws.Range("D42").Value = "'6.617"
$ws.Range("E42").Value = "  -0.55%  "
$ws.Range("D43").Value = "'116.68"
$ws.Range("E43").Value = "  -0.42%  "
$ws.Range("D44").Value = "'0.1519"
$ws.Range("E44").Value = "  +0.34%  "
$ws.Range("D45").Value = "'0.4915"
$ws.Range("E45").Value = "  +0.76%  "
$ws.Range("D46").Value = "'10.20"
$ws.Range("E46").Value = "  +0.81%  "
$ws.Range("D47").Value = "'1.002"
$ws.Range("E47").Value = "  +0.20%  "
$ws.Range("D48").Value = "'1.640"
$ws.Range("E48").Value = "  +1.46%  "
$ws.Range("D49").Value = "'38.55"
$ws.Range("E49").Value = "  +2.84%  "
$ws.Range("D50").Value = "'64.27"
$ws.Range("E50").Value = "  +0.50%  "
$ws.Range("D51").Value = "'0.05961"
$ws.Range("E51").Value = "  +0.37%  "
